$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet/tab to reflect the new "through" date
$ws.Name = "Through 2022-04-22"

# Update the header label in I1 (shared string "2022 (through 04-21)" -> "2022 (through 04-22)")
$ws.Range("I1").Value = "2022 (through 04-22)"

# Update April's year-to-date count
$ws.Range("I5").Value = 92

# Update the Total row's year-to-date count
$ws.Range("I14").Value = 528
